# fix: hilangkan referral di mou
#
# Removes the three "referral" placeholder paragraphs (the
# "${referral}{{", "${referral_signature:100px:100px:ratio=true}" and
# "${referral_name}email}}" paragraphs, including their paragraph
# marks) from the MOU template, leaving the surrounding empty
# paragraphs untouched.

$d = $word.ActiveDocument

# Locate the first paragraph to remove: it starts with the literal
# text "${referral}{{".
$startRng = $d.Content.Duplicate
$startRng.Find.MatchWildcards = $false
$startRng.Find.Execute("`${referral}{{") | Out-Null
$startRng.Expand(4) | Out-Null   # wdParagraph -> whole paragraph incl. mark

# Locate the last paragraph to remove: it ends with the literal text
# "email}}" (the tail of "${referral_name}email}}").
$endRng = $d.Content.Duplicate
$endRng.Find.MatchWildcards = $false
$endRng.Find.Execute("email}}") | Out-Null
$endRng.Expand(4) | Out-Null     # wdParagraph -> whole paragraph incl. mark

# Delete everything from the start of the first paragraph through the
# end (incl. paragraph mark) of the last paragraph - this removes all
# three referral paragraphs in one shot and merges the surrounding
# empty paragraphs back together exactly as before.
$delRng = $d.Range($startRng.Start, $endRng.End)
$delRng.Delete()
